# Add data from 3/27/2020 (new rows 457-530) to the houstonNumbers sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: new rows use the same date-formatted style as the existing
#     date column (row 456 already carries the desired numFmt style), so
#     copy that cell's formatting into the new range before writing values.
$dateStyleSource = $ws.Cells.Item(456, 1)
$newDateRange = $ws.Range("A457:A530")
$dateStyleSource.Copy($newDateRange)

# All the new rows are dated 3/27/2020 (Excel serial date 43917).
$ws.Range("A457:A530").Value = 43917

# --- Column B (Sex) ---
$ws.Range("B472:B479").Value = "F"
$ws.Range("B480:B484").Value = "M"
$ws.Range("B530").Value = "F"

# --- Column C (Age) ---
$ws.Range("C530").Value = "50-60"

# --- Column D (County) ---
$ws.Range("D457:D471").Value = "Fort Bend"
$ws.Range("D472:D484").Value = "Brazoria"
$ws.Range("D485:D493").Value = "Galveston"
$ws.Range("D494").Value = "Matagorda"
$ws.Range("D495:D519").Value = "Harris"
$ws.Range("D520").Value = "Houston"
$ws.Range("D521:D526").Value = "Montgomery"
$ws.Range("D527:D529").Value = "Brazos"
$ws.Range("D530").Value = "Matagorda"

# Match the author's final view/selection state from the saved workbook.
$ws.Range("C511").Select()
